# "Casos pendientes al 95%"
#
# The tracker lost its "James Andres Urquiza" assignment; that row is
# reassigned to "Luis Carlos Rincon Gordo" (the workbook owner), whose
# e-mail address is now linked both on his new row and on Frank Stiven's
# row (which previously had no e-mail hyperlink at all).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- text updates -----------------------------------------------------
# Set B2's value before touching row 4 so that the shared-strings table
# gets the e-mail string introduced ahead of the new person's name,
# matching the order used by the authored workbook.
$ws.Range("B2").Value = "Luis.RinconG@axity.com"

# Row 4: replace the old assignee with the new one, and refresh the
# e-mail text next to it.
$ws.Range("A4").Value = "Luis Carlos Rincon Gordo"
$ws.Range("B4").Value = "Luis.RinconG@axity.com"

# --- hyperlinks ---------------------------------------------------------
# Create the B4 hyperlink first so it becomes rId1, then B2 as rId2.
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:Luis.RinconG@axity.com", [Type]::Missing, [Type]::Missing, "Luis.RinconG@axity.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Luis.RinconG@axity.com", [Type]::Missing, [Type]::Missing, "Luis.RinconG@axity.com") | Out-Null

# Adding a hyperlink resets the cell's look; reapply the workbook's
# existing built-in "Hyperlink" cell style so no new style is needed.
$ws.Range("B4").Style = "Hipervínculo"
$ws.Range("B2").Style = "Hipervínculo"

# --- stray formatted cell ------------------------------------------------
# A leftover formatted (underlined) empty cell shows up at D7, matching
# the style already used by the empty C4 cell.
$ws.Range("D7").Font.Underline = $true

# --- selection ------------------------------------------------------------
$ws.Range("D7").Select() | Out-Null

Write-Host "Edit complete"
